$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.307.13"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.50%  "
$ws.Range("D3").Value = "'3.426.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.18%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'413.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.34%  "
$ws.Range("E6").Value = "  -1.09%  "
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  -2.02%  "
$ws.Range("E10").Value = "  -0.98%  "
$ws.Range("D11").Value = "'42.77"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.79%  "
$ws.Range("E12").Value = "  +3.26%  "
$ws.Range("D13").Value = "'9.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "'3.964.85"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E16").Value = "  -3.56%  "
$ws.Range("D17").Value = "'3.403.13"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'12.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.32%  "
$ws.Range("E19").Value = "  -1.24%  "
$ws.Range("D20").Value = "'62.298.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("D21").Value = "'473.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  +2.58%  "
$ws.Range("D24").Value = "'13.11"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.11%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'9.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'33.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.68%  "
$ws.Range("D28").Value = "'4.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.73%  "
$ws.Range("D29").Value = "'7.72"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").Value = "'11.89"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.19%  "
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D34").Value = "'40.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.46%  "
$ws.Range("D36").Value = "'57.94"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.34%  "
$ws.Range("E37").Value = "  -2.75%  "
$ws.Range("E38").Value = "  +0.05%  "
$ws.Range("D39").Value = "'3.01"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").Value = "'145.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.49%  "
$ws.Range("D44").Value = "'2.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.73%  "
$ws.Range("D45").Value = "'4.34"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.40%  "
$ws.Range("E46").Value = "  +3.56%  "
$ws.Range("D47").Value = "'2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +18.86%  "
$ws.Range("E48").Value = "  -2.42%  "
$ws.Range("D49").Value = "'0.0₃0542"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +25.33%  "
$ws.Range("D50").Value = "'22.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'112.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.28%  "
